# Fix bug: a number of cells in columns F,G,H,I,K,L,M,N (the "conti-ord"
# and "contin-nom" test columns) were incorrectly recorded as 2 and must
# be corrected to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellsToFix = @(
    "F2","H2","I2","L2","M2",
    "F3","H3","I3","L3","M3",
    "F4","H4","I4","L4","M4",
    "F5","H5","I5","L5","M5",
    "L6","M6",
    "L7","M7",
    "L10","M10",
    "K11","M11",
    "F13","H13","L13","M13",
    "F14","H14","I14","L14","M14",
    "F15","H15","I15","L15","M15",
    "F17","H17","L17","M17",
    "L18","M18",
    "F20","H20",
    "F21","H21",
    "F22","H22",
    "F23","H23","L23","M23",
    "F24","H24","L24","M24",
    "F25","H25","I25","L25","M25",
    "F26","H26",
    "F27","H27","I27","L27","M27",
    "F28","H28","I28","L28","M28",
    "F29","H29","I29","L29","M29","N29",
    "F30","H30","I30",
    "F31","H31","I31","L31","M31",
    "F32","H32","I32","L32","M32",
    "F34","H34",
    "F35","H35","L35","M35",
    "F36","H36","I36","L36","M36",
    "F37","H37","L37","M37",
    "L38","M38",
    "F39","H39","I39",
    "L40","M40","N40",
    "F41","H41",
    "F42","H42","I42",
    "F43","H43","K43","L43",
    "F46","H46",
    "F48","H48","L48","M48",
    "F49","H49","I49",
    "F51","H51","L51","M51",
    "F52","H52","I52","L52","M52","N52",
    "F53","I53",
    "F54","H54","L54","M54",
    "L55","M55",
    "F56","H56",
    "L59","M59",
    "F64","H64","I64","L64","M64",
    "L65","M65",
    "F68","H68","I68","L68","M68",
    "F69","H69","I69",
    "F75","H75",
    "F76","G76","H76","I76",
    "L77","M77",
    "L78","M78"
)

foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = 1
}
